$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 10010802
$ws.Range("I74").Value = 16671001
$ws.Range("K74").Value = 16671001
$ws.Range("M74").Value = -16670065

# Row 77
$ws.Range("H77").Value = 10010802
$ws.Range("I77").Value = 16671001
$ws.Range("K77").Value = 83355005
$ws.Range("M77").Value = -83350325

# Row 98
$ws.Range("H98").Value = 9812.6875
$ws.Range("I98").Value = 8667.333000000001
$ws.Range("J98").Value = 10499.9
$ws.Range("K98").Value = 8667.333000000001
$ws.Range("L98").Value = 10499.9
$ws.Range("M98").Value = -7169.333000000001
$ws.Range("N98").Value = -13495.9

# Row 113
$ws.Range("H113").Value = 4738.077
$ws.Range("I113").Value = 2173.5715
$ws.Range("J113").Value = 7730
$ws.Range("K113").Value = 2173.5715
$ws.Range("L113").Value = 7730
$ws.Range("M113").Value = 1080.4285
$ws.Range("N113").Value = -14238

# Row 116
$ws.Range("H116").Value = 506032.84
$ws.Range("I116").Value = 2501948.8
$ws.Range("J116").Value = 7053.875
$ws.Range("K116").Value = 2501948.8
$ws.Range("L116").Value = 7053.875
$ws.Range("M116").Value = -2498506.8
$ws.Range("N116").Value = -13937.875

# Row 122
$ws.Range("H122").Value = 9812.6875
$ws.Range("I122").Value = 8667.333000000001
$ws.Range("J122").Value = 10499.9
$ws.Range("K122").Value = 26001.999
$ws.Range("L122").Value = 31499.7
$ws.Range("M122").Value = -23551.999
$ws.Range("N122").Value = -36399.7

# Row 129
$ws.Range("H129").Value = 813.1031
$ws.Range("I129").Value = 285.1111
$ws.Range("J129").Value = 867.1023
$ws.Range("K129").Value = 855.3333
$ws.Range("L129").Value = 2601.3069
$ws.Range("M129").Value = 4144.6667
$ws.Range("N129").Value = -12601.3069

# Row 132
$ws.Range("H132").Value = 45457252
$ws.Range("I132").Value = 50002130
$ws.Range("J132").Value = 8500
$ws.Range("K132").Value = 150006390
$ws.Range("L132").Value = 25500
$ws.Range("M132").Value = -150003860
$ws.Range("N132").Value = -30560

# Row 133
$ws.Range("H133").Value = 67173.336
$ws.Range("J133").Value = 67173.336
$ws.Range("L133").Value = 67173.336
$ws.Range("N133").Value = -77293.336

$ws = $wb.Worksheets.Item("ARM")
# Row 101
$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490

# Row 134
$ws.Range("H134").Value = 39998
$ws.Range("J134").Value = 39998
$ws.Range("L134").Value = 39998
$ws.Range("N134").Value = -50138

# Row 137
$ws.Range("H137").Value = 40953.332
$ws.Range("J137").Value = 40953.332
$ws.Range("L137").Value = 40953.332
$ws.Range("N137").Value = -51153.332

$ws = $wb.Worksheets.Item("BSM")
# Row 137
$ws.Range("H137").Value = 50760
$ws.Range("J137").Value = 50760
$ws.Range("L137").Value = 50760
$ws.Range("N137").Value = -60960

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2985.75
$ws.Range("I58").Value = 1801.5103
$ws.Range("K58").Value = 1801.5103
$ws.Range("M58").Value = -1598.5103

# Row 132
$ws.Range("H132").Value = 2696.7693
$ws.Range("I132").Value = 2211.3125
$ws.Range("K132").Value = 6633.9375
$ws.Range("M132").Value = -4103.9375

# Row 134
$ws.Range("H134").Value = 1736.5
$ws.Range("I134").Value = 1107.7693
$ws.Range("K134").Value = 3323.3079
$ws.Range("M134").Value = -788.3078999999998

# Row 136
$ws.Range("H136").Value = 2985.75
$ws.Range("I136").Value = 1801.5103
$ws.Range("K136").Value = 5404.5309
$ws.Range("M136").Value = -2854.5309

# Row 140
$ws.Range("H140").Value = 92196.37
$ws.Range("J140").Value = 92196.37
$ws.Range("L140").Value = 92196.37
$ws.Range("N140").Value = -102556.37

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 477880.72
$ws.Range("I5").Value = 500.86667
$ws.Range("J5").Value = 1028703.6
$ws.Range("K5").Value = 1502.60001
$ws.Range("L5").Value = 3086110.8
$ws.Range("M5").Value = -1390.60001
$ws.Range("N5").Value = -3086334.8

# Row 122
$ws.Range("H122").Value = 3067.2888
$ws.Range("J122").Value = 3508.8647
$ws.Range("L122").Value = 31579.7823
$ws.Range("N122").Value = -36479.7823

# Row 132
$ws.Range("H132").Value = 2355.7742
$ws.Range("I132").Value = 945.3333
$ws.Range("J132").Value = 3246.5789
$ws.Range("K132").Value = 8507.9997
$ws.Range("L132").Value = 29219.2101
$ws.Range("M132").Value = -5977.9997
$ws.Range("N132").Value = -34279.2101

# Row 135
$ws.Range("H135").Value = 477880.72
$ws.Range("I135").Value = 500.86667
$ws.Range("J135").Value = 1028703.6
$ws.Range("K135").Value = 4507.80003
$ws.Range("L135").Value = 9258332.4
$ws.Range("M135").Value = -1972.80003
$ws.Range("N135").Value = -9263402.4

# Row 137
$ws.Range("H137").Value = 846
$ws.Range("I137").Value = 846
$ws.Range("K137").Value = 2538
$ws.Range("M137").Value = 2562

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -708

# Row 37
$ws.Range("H37").Value = 1000
$ws.Range("I37").Value = 1000
$ws.Range("K37").Value = 1000
$ws.Range("M37").Value = -723

# Row 46
$ws.Range("H46").Value = 31250.666
$ws.Range("J46").Value = 31250.666
$ws.Range("L46").Value = 31250.666
$ws.Range("N46").Value = -31562.666

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 126
$ws.Range("H126").Value = 3322.14
$ws.Range("I126").Value = 2931.081
$ws.Range("J126").Value = 4435.154
$ws.Range("K126").Value = 8793.243
$ws.Range("L126").Value = 13305.462
$ws.Range("M126").Value = -6323.243
$ws.Range("N126").Value = -18245.462

# Row 137
$ws.Range("H137").Value = 40367.2
$ws.Range("J137").Value = 40367.2
$ws.Range("L137").Value = 40367.2
$ws.Range("N137").Value = -50567.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7628
$ws.Range("I7").Value = 7582.857
$ws.Range("K7").Value = 7582.857
$ws.Range("M7").Value = -7470.857

# Row 29
$ws.Range("H29").Value = 7999.5
$ws.Range("I29").Value = 6000
$ws.Range("J29").Value = 9999
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 9999
$ws.Range("M29").Value = -5705
$ws.Range("N29").Value = -10589

# Row 40
$ws.Range("H40").Value = 4771.2856
$ws.Range("I40").Value = 3066.5
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 3066.5
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -2930.5
$ws.Range("N40").Value = -15272

# Row 126
$ws.Range("H126").Value = 7628
$ws.Range("I126").Value = 7582.857
$ws.Range("K126").Value = 22748.571
$ws.Range("M126").Value = -20278.571

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3779.7585
$ws.Range("I122").Value = 2391.5
$ws.Range("J122").Value = 8142.857
$ws.Range("K122").Value = 7174.5
$ws.Range("L122").Value = 24428.571
$ws.Range("M122").Value = -4724.5
$ws.Range("N122").Value = -29328.571

# Row 126
$ws.Range("H126").Value = 4466.857
$ws.Range("I126").Value = 2058
$ws.Range("J126").Value = 10489
$ws.Range("K126").Value = 6174
$ws.Range("L126").Value = 31467
$ws.Range("M126").Value = -3704
$ws.Range("N126").Value = -36407

# Row 136
$ws.Range("H136").Value = 14176.286
$ws.Range("I136").Value = 20897
$ws.Range("J136").Value = 10442.556
$ws.Range("K136").Value = 62691
$ws.Range("L136").Value = 31327.668
$ws.Range("M136").Value = -60141
$ws.Range("N136").Value = -36427.66800000001
